$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 25.50000000000055
$ws.Range("G2").Value = [double]"2.149391775674303e-13"
$ws.Range("H2").Value = [double]"3.21241161012888e-12"
$ws.Range("I2").Value = 0.31455062890494
$ws.Range("K2").Value = 39.06952453244617
$ws.Range("L2").Value = "[28.35602567676036, 49.78302338813199]"
$ws.Range("M2").Value = [double]"1.435740415445252e-11"
$ws.Range("N2").Value = [double]"2.871480830890505e-11"
$ws.Range("O2").Value = 1.389973926813502
$ws.Range("P2").Value = "[1.0755001877154244, 1.7044476659115793]"
$ws.Range("Q2").Value = [double]"1.554312234475219e-15"
$ws.Range("R2").Value = [double]"3.108624468950438e-15"
$ws.Range("S2").Value = 56.84293922338828
$ws.Range("T2").Value = "[50.29717234672082, 63.388706100055735]"
$ws.Range("W2").Value = 19.85885885885929
$ws.Range("X2").Value = 18.58258258258298
$ws.Range("Y2").Value = 21.13513513513559

# --- Row 3 ---
$ws.Range("E3").Value = 23.62000000000025
$ws.Range("G3").Value = [double]"5.433292260548228e-09"
$ws.Range("H3").Value = [double]"1.486545953412485e-08"
$ws.Range("K3").Value = 40.05319326779832
$ws.Range("L3").Value = "[23.785787600284067, 56.320598935312574]"
$ws.Range("M3").Value = [double]"2.988230624723087e-06"
$ws.Range("N3").Value = [double]"2.988230624723087e-06"
$ws.Range("O3").Value = -0.8050527720910781
$ws.Range("P3").Value = "[-1.2075791581366175, -0.4025263860455386]"
$ws.Range("Q3").Value = [double]"0.0001212246300852371"
$ws.Range("R3").Value = [double]"0.0001212246300852371"
$ws.Range("S3").Value = 64.17670327655622
$ws.Range("T3").Value = "[55.65597351141869, 72.69743304169376]"
$ws.Range("W3").Value = 3.026386386386417
$ws.Range("X3").Value = 1.513193193193206
$ws.Range("Y3").Value = 4.539579579579628
